# Weather log sheet update:
#  - header row gains a "location" label (replacing the old "Temperature"
#    text in B1) and extends from J (Temperature8) through T (Temperature18)
#  - the old gilgit/skardu/srinagar block (A2:G4) is replaced by a new
#    dusseldorf/Nice/Marseille/Monte Carlo block (A2:A5) whose scraped
#    "temperature" readings landed (mis-scraped, ChromeDriver session-id
#    text and all) in columns B/C, with clean numeric duplicates in R/S.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- header row -------------------------------------------------------
$ws.Range("B1").Value() = "location"

$ws.Range("K1").Value() = "Temperature9"
$ws.Range("L1").Value() = "Temperature10"
$ws.Range("M1").Value() = "Temperature11"
$ws.Range("N1").Value() = "Temperature12"
$ws.Range("O1").Value() = "Temperature13"
$ws.Range("P1").Value() = "Temperature14"
$ws.Range("Q1").Value() = "Temperature15"
$ws.Range("R1").Value() = "Temperature16"
$ws.Range("S1").Value() = "Temperature17"
$ws.Range("T1").Value() = "Temperature18"

# new header cells need the same yellow header fill as the rest of row 1
$ws.Range("K1:T1").Interior.Color = 65535

# --- clear the old gilgit/skardu/srinagar reading columns -------------
$ws.Range("D2:G4").ClearContents()

# --- new places ---------------------------------------------------------
$ws.Range("A2").Value() = "dusseldorf"
$ws.Range("A3").Value() = "Nice"
$ws.Range("A4").Value() = "Marseille"
$ws.Range("A5").Value() = "Monte Carlo"

# --- scraped readings (raw ChromeDriver artifacts in B/C) -------------
$ws.Range("B2").Value() = '16.04.$[[ChromeDriver: chrome on XP (70d13a5d04fcbe4353588bb598b793d3)] -> id: xPat]'
$ws.Range("C2").Value() = '16.04.$[[ChromeDriver: chrome on XP (70d13a5d04fcbe4353588bb598b793d3)] -> id: xPat]'

$ws.Range("B3").Value() = '13.58.$[[ChromeDriver: chrome on XP (6467ec0d74fea02c0a3ac09ad57bcefb)] -> id: xPat]'
$ws.Range("C3").Value() = '13.58.$[[ChromeDriver: chrome on XP (6467ec0d74fea02c0a3ac09ad57bcefb)] -> id: xPat]'

$ws.Range("B4").Value() = '20.36.$[[ChromeDriver: chrome on XP (d2699a01bbe2d02d9fdb2c4a0859a00f)] -> id: xPat]'
$ws.Range("C4").Value() = '20.36.$[[ChromeDriver: chrome on XP (d2699a01bbe2d02d9fdb2c4a0859a00f)] -> id: xPat]'

$ws.Range("B5").Value() = '21.96.$[[ChromeDriver: chrome on XP (52ab270f72e9e43db6022c5f290b4b9d)] -> id: xPat]'
$ws.Range("C5").Value() = '21.96.$[[ChromeDriver: chrome on XP (52ab270f72e9e43db6022c5f290b4b9d)] -> id: xPat]'

# --- cleaned-up numeric readings (R/S) ---------------------------------
$ws.Range("R2").Value() = "13.3."
$ws.Range("S2").Value() = "13.3."

$ws.Range("R3").Value() = "13.28."
$ws.Range("S3").Value() = "13.28."

$ws.Range("R4").Value() = "19.12."
$ws.Range("S4").Value() = "19.12."

$ws.Range("R5").Value() = "21.81."
$ws.Range("S5").Value() = "21.81."

# --- column widths (best-fit-ish for the new wide text columns) -------
$ws.Range("B1:C1").ColumnWidth = 83.03125
$ws.Range("D1:K1").ColumnWidth = 13.85546875
$ws.Range("L1:S1").ColumnWidth = 15

# --- final selection ----------------------------------------------------
$ws.Range("B1").Select()
